# Fix the "QR" calculated column in Tabela1: the old formula used INT(B2)
# which fails (#VALUE!) because column B holds text like "SP001". Replace it
# with VALUE(SUBSTITUTE(B2,"SP","")) so the spool number is parsed correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo  = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item(16)   # "QR" column

$col.DataBodyRange.Formula = '=CONCATENATE("2$-", A2, "$ID-", VALUE(SUBSTITUTE(B2, "SP", "")))'

# Column P now shows longer text results instead of #VALUE!, so widen it
# to fit the new content (matches the workbook's resulting column width).
$ws.Columns.Item(16).ColumnWidth = 26.14

# Move/collapse the active selection onto P2 (single cell), as left by the
# author after making the edit.
$ws.Range("P2").Select() | Out-Null
